$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

# Row 112: clear the empty H112 inline-string cell (Actions column)
$ws.Range("H112").Value = $null

# Row 113: add reference QB#5079 in column D
$ws.Range("D113").Value = "QB#5079"

# New row 116
$ws.Range("A116").Value = "خبراء الفلاتر"
$ws.Range("B116").Value = "'2025-07-28"
$ws.Range("B116").ClearFormats()
$ws.Range("C116").Value = "#1: استيكر منتج - طباعة ديجيتال - مقاس 10*15 سم | Qty: 100 | Price: 0 | Total: 0 | VAT: 0"
$ws.Range("E116").Value = 0
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("I116").Value = $false

# New row 117
$ws.Range("A117").Value = "شاي فال"
$ws.Range("B117").Value = "'2025-07-28"
$ws.Range("B117").ClearFormats()
$ws.Range("C117").Value = "#1: استيكرات النظافة - بلاش على قولة سامي | Qty: 5 | Price: 0 | Total: 0 | VAT: 0"
$ws.Range("E117").Value = 0
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("I117").Value = $false

# New row 118
$ws.Range("A118").Value = "مجمع قبل وبعد"
$ws.Range("B118").Value = "'2025-07-28"
$ws.Range("B118").ClearFormats()
$ws.Range("C118").Value = "#1: بطاقات آي دي موظفين | Qty: 3 | Price: 28 | Total: 84 | VAT: 12.6"
$ws.Range("E118").Value = 84
$ws.Range("F118").Value = 12.6
$ws.Range("G118").Value = 96.59999999999999
$ws.Range("I118").Value = $false

# New row 119
$ws.Range("A119").Value = "روافد القهوة"
$ws.Range("B119").Value = "'2025-07-29"
$ws.Range("B119").ClearFormats()
$ws.Range("C119").Value = "#1: استيكر مقاس 24*13 سم - طباعة ديجيتال - مع سلوفان مطفي | Qty: 150 | Price: 1 | Total: 150 | VAT: 22.5"
# D119 is present but blank (empty text), mirrors the pre-existing blank H112 cell
$ws.Range("D119").Value = "'"
$ws.Range("D119").ClearFormats()
$ws.Range("E119").Value = 150
$ws.Range("F119").Value = 22.5
$ws.Range("G119").Value = 172.5
# H119 is present but blank (empty text) too
$ws.Range("H119").Value = "'"
$ws.Range("H119").ClearFormats()
$ws.Range("I119").Value = $false
